$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.934.12"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "2.217.06"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'262.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").Value = "'87.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.82%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("D10").Value = "'45.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.95%  "
$ws.Range("D11").Value = "'0.0922"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("D12").Value = "'7.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.24%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").Value = "2.548.30"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "'14.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "2.209.11"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "'0.787"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "43.856.45"
$ws.Range("E18").Value = "  +2.34%  "
$ws.Range("D19").Value = "'0.0000104"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "'5.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "'70.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").Value = "'2.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.48%  "
$ws.Range("D23").Value = "'232.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").Value = "'8.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.61%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'10.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").Value = "'3.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.72%  "
$ws.Range("D28").Value = "'39.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.74%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.36%  "
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").Value = "'174.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "'20.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Value = "'0.0880"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "'5.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.01%  "
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("E36").Value = "  +4.38%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0359"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.39%  "
$ws.Range("D39").Value = "'3.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.00%  "
$ws.Range("D40").Value = "'12.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("D41").Value = "'65.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.21%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").Value = "'5.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.87%  "
$ws.Range("D44").Value = "'0.202"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").Value = "'101.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'8.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0983"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "'1.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.51%  "
$ws.Range("D50").Value = "'0.447"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("D51").Value = "'1.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.22%  "
